$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3498.5
$ws.Range("J32").Value = 3498.5
$ws.Range("L32").Value = 3498.5
$ws.Range("N32").Value = -4150.5
$ws.Range("H96").Value = 947.7273
$ws.Range("I96").Value = 756.6667
$ws.Range("K96").Value = 2270.0001
$ws.Range("M96").Value = -897.0001000000002
$ws.Range("H107").Value = 579.7143
$ws.Range("I107").Value = 603.53845
$ws.Range("K107").Value = 603.53845
$ws.Range("M107").Value = 1316.46155
$ws.Range("H112").Value = 6409.8823
$ws.Range("J112").Value = 6717.8667
$ws.Range("L112").Value = 20153.6001
$ws.Range("N112").Value = -22369.6001
$ws.Range("H138").Value = 3174.1638
$ws.Range("J138").Value = 2652.9111
$ws.Range("L138").Value = 7958.7333
$ws.Range("N138").Value = -18238.7333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3115
$ws.Range("I45").Value = 2827.1428
$ws.Range("J45").Value = 3518
$ws.Range("K45").Value = 2827.1428
$ws.Range("L45").Value = 3518
$ws.Range("M45").Value = -2450.1428
$ws.Range("N45").Value = -4272
$ws.Range("H61").Value = 3503.9119
$ws.Range("I61").Value = 3497.75
$ws.Range("J61").Value = 3518.7
$ws.Range("K61").Value = 3497.75
$ws.Range("L61").Value = 3518.7
$ws.Range("M61").Value = -3285.75
$ws.Range("N61").Value = -3942.7
$ws.Range("H74").Value = 2091.7778
$ws.Range("I74").Value = 1964.2632
$ws.Range("J74").Value = 2394.625
$ws.Range("K74").Value = 1964.2632
$ws.Range("L74").Value = 2394.625
$ws.Range("M74").Value = -1090.2632
$ws.Range("N74").Value = -4142.625
$ws.Range("H77").Value = 2091.7778
$ws.Range("I77").Value = 1964.2632
$ws.Range("J77").Value = 2394.625
$ws.Range("K77").Value = 9821.316000000001
$ws.Range("L77").Value = 11973.125
$ws.Range("M77").Value = -5453.316000000001
$ws.Range("N77").Value = -20709.125
$ws.Range("H88").Value = 9832.25
$ws.Range("I88").Value = 1291.7142
$ws.Range("J88").Value = 21789
$ws.Range("K88").Value = 1291.7142
$ws.Range("L88").Value = 21789
$ws.Range("M88").Value = -885.7141999999999
$ws.Range("N88").Value = -22601
$ws.Range("H91").Value = 9832.25
$ws.Range("I91").Value = 1291.7142
$ws.Range("J91").Value = 21789
$ws.Range("K91").Value = 1291.7142
$ws.Range("L91").Value = 21789
$ws.Range("M91").Value = 112.2858000000001
$ws.Range("N91").Value = -24597
$ws.Range("H111").Value = 90000
$ws.Range("J111").Value = 90000
$ws.Range("L111").Value = 90000
$ws.Range("N111").Value = -98180
$ws.Range("H132").Value = 2263.3914
$ws.Range("I132").Value = 2037.95
$ws.Range("K132").Value = 6113.85
$ws.Range("M132").Value = -3583.85
$ws.Range("H136").Value = 3503.9119
$ws.Range("I136").Value = 3497.75
$ws.Range("J136").Value = 3518.7
$ws.Range("K136").Value = 10493.25
$ws.Range("L136").Value = 10556.1
$ws.Range("M136").Value = -7943.25
$ws.Range("N136").Value = -15656.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 94998.5
$ws.Range("J55").Value = 94998.5
$ws.Range("L55").Value = 94998.5
$ws.Range("N55").Value = -95544.5
$ws.Range("H86").Value = 8587.637000000001
$ws.Range("I86").Value = 3581.5334
$ws.Range("K86").Value = 3581.5334
$ws.Range("M86").Value = -2458.5334
$ws.Range("H89").Value = 8587.637000000001
$ws.Range("I89").Value = 3581.5334
$ws.Range("K89").Value = 17907.667
$ws.Range("M89").Value = -12291.667
$ws.Range("H94").Value = 24252.092
$ws.Range("I94").Value = 10979.167
$ws.Range("K94").Value = 10979.167
$ws.Range("M94").Value = -10528.167
$ws.Range("H105").Value = 3289.2122
$ws.Range("I105").Value = 2991.4614
$ws.Range("K105").Value = 2991.4614
$ws.Range("M105").Value = -1244.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 8012.5293
$ws.Range("I86").Value = 7640.2
$ws.Range("J86").Value = 8544.429
$ws.Range("K86").Value = 7640.2
$ws.Range("L86").Value = 8544.429
$ws.Range("M86").Value = -6517.2
$ws.Range("N86").Value = -10790.429
$ws.Range("H89").Value = 8012.5293
$ws.Range("I89").Value = 7640.2
$ws.Range("J89").Value = 8544.429
$ws.Range("K89").Value = 38201
$ws.Range("L89").Value = 42722.145
$ws.Range("M89").Value = -32585
$ws.Range("N89").Value = -53954.145
$ws.Range("H99").Value = 1666.9286
$ws.Range("I99").Value = 1687.4615
$ws.Range("K99").Value = 1687.4615
$ws.Range("M99").Value = -189.4614999999999
$ws.Range("H107").Value = 795.6316
$ws.Range("I107").Value = 716.8333
$ws.Range("J107").Value = 930.7143
$ws.Range("K107").Value = 716.8333
$ws.Range("L107").Value = 930.7143
$ws.Range("M107").Value = 1203.1667
$ws.Range("N107").Value = -4770.7143
$ws.Range("H126").Value = 1666.9286
$ws.Range("I126").Value = 1687.4615
$ws.Range("K126").Value = 5062.3845
$ws.Range("M126").Value = -2592.3845
$ws.Range("H132").Value = 2662.2163
$ws.Range("I132").Value = 2458.6875
$ws.Range("J132").Value = 3964.8
$ws.Range("K132").Value = 7376.0625
$ws.Range("L132").Value = 11894.4
$ws.Range("M132").Value = -4846.0625
$ws.Range("N132").Value = -16954.4
$ws.Range("H134").Value = 1125.8182
$ws.Range("J134").Value = 1218.1666
$ws.Range("L134").Value = 3654.4998
$ws.Range("N134").Value = -8724.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 584.4
$ws.Range("I5").Value = 543.25
$ws.Range("J5").Value = 749
$ws.Range("K5").Value = 1629.75
$ws.Range("L5").Value = 2247
$ws.Range("M5").Value = -1517.75
$ws.Range("N5").Value = -2471
$ws.Range("H38").Value = 180.23529
$ws.Range("I38").Value = 81.44444
$ws.Range("J38").Value = 291.375
$ws.Range("K38").Value = 244.33332
$ws.Range("L38").Value = 874.125
$ws.Range("M38").Value = 102.66668
$ws.Range("N38").Value = -1568.125
$ws.Range("H64").Value = 6531
$ws.Range("I64").Value = 1750
$ws.Range("J64").Value = 7487.2
$ws.Range("K64").Value = 5250
$ws.Range("L64").Value = 22461.6
$ws.Range("M64").Value = -4980
$ws.Range("N64").Value = -23001.6
$ws.Range("H67").Value = 6531
$ws.Range("I67").Value = 1750
$ws.Range("J67").Value = 7487.2
$ws.Range("K67").Value = 5250
$ws.Range("L67").Value = 22461.6
$ws.Range("M67").Value = -4314
$ws.Range("N67").Value = -24333.6
$ws.Range("H97").Value = 399.5
$ws.Range("I97").Value = 299
$ws.Range("K97").Value = 897
$ws.Range("M97").Value = -401
$ws.Range("H133").Value = 12494.618
$ws.Range("J133").Value = 13223.1
$ws.Range("L133").Value = 39669.3
$ws.Range("N133").Value = -49789.3
$ws.Range("H135").Value = 584.4
$ws.Range("I135").Value = 543.25
$ws.Range("J135").Value = 749
$ws.Range("K135").Value = 4889.25
$ws.Range("L135").Value = 6741
$ws.Range("M135").Value = -2354.25
$ws.Range("N135").Value = -11811
$ws.Range("H136").Value = 6091
$ws.Range("I136").Value = 3227.889
$ws.Range("K136").Value = 9683.667000000001
$ws.Range("M136").Value = -4583.667000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2544.2354
$ws.Range("I80").Value = 2489.7334
$ws.Range("K80").Value = 2489.7334
$ws.Range("M80").Value = -1491.7334
$ws.Range("H83").Value = 2544.2354
$ws.Range("I83").Value = 2489.7334
$ws.Range("K83").Value = 12448.667
$ws.Range("M83").Value = -7456.667000000001
$ws.Range("H97").Value = 1055.7368
$ws.Range("J97").Value = 3173.6667
$ws.Range("L97").Value = 3173.6667
$ws.Range("N97").Value = -4165.6667
$ws.Range("H98").Value = 16518.25
$ws.Range("J98").Value = 16518.25
$ws.Range("L98").Value = 16518.25
$ws.Range("N98").Value = -22508.25
$ws.Range("H132").Value = 1797.0488
$ws.Range("I132").Value = 1794.1842
$ws.Range("J132").Value = 1833.3334
$ws.Range("K132").Value = 5382.5526
$ws.Range("L132").Value = 5500.0002
$ws.Range("M132").Value = -2852.5526
$ws.Range("N132").Value = -10560.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2436.25
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 2581.6667
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 2581.6667
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -2805.6667
$ws.Range("H82").Value = 6404.75
$ws.Range("J82").Value = 6969.7144
$ws.Range("L82").Value = 6969.7144
$ws.Range("N82").Value = -7691.7144
$ws.Range("H85").Value = 6404.75
$ws.Range("J85").Value = 6969.7144
$ws.Range("L85").Value = 6969.7144
$ws.Range("N85").Value = -9465.714400000001
$ws.Range("H108").Value = 79999.5
$ws.Range("J108").Value = 79999.5
$ws.Range("L108").Value = 79999.5
$ws.Range("N108").Value = -87679.5
$ws.Range("H109").Value = 64991.5
$ws.Range("J109").Value = 64991.5
$ws.Range("L109").Value = 64991.5
$ws.Range("N109").Value = -67765.5
$ws.Range("H126").Value = 2436.25
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2581.6667
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 7745.000100000001
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -12685.0001
$ws.Range("H132").Value = 2271.1538
$ws.Range("I132").Value = 1792.6279
$ws.Range("J132").Value = 3206.4546
$ws.Range("K132").Value = 5377.8837
$ws.Range("L132").Value = 9619.363799999999
$ws.Range("M132").Value = -2847.8837
$ws.Range("N132").Value = -14679.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 7355.1333
$ws.Range("I107").Value = 10264.7
$ws.Range("K107").Value = 30794.1
$ws.Range("M107").Value = -28874.1
$ws.Range("H109").Value = 80000
$ws.Range("J109").Value = 80000
$ws.Range("L109").Value = 80000
$ws.Range("N109").Value = -82774
$ws.Range("H136").Value = 3081.963
$ws.Range("I136").Value = 2429.2173
$ws.Range("J136").Value = 6835.25
$ws.Range("K136").Value = 7287.651899999999
$ws.Range("L136").Value = 20505.75
$ws.Range("M136").Value = -4737.651899999999
$ws.Range("N136").Value = -25605.75
